$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.492.88'
$ws.Range('E2').Value = '  -7.72%  '
$ws.Range('D3').Value = '2.547.29'
$ws.Range('E3').Value = '  -1.94%  '
$ws.Range('E4').Value = '  +0.06%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '295.42'
$ws.Range('E5').Value = '  -5.19%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '90.70'
$ws.Range('E6').Value = '  -8.22%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.570'
$ws.Range('E7').Value = '  -4.40%  '
$ws.Range('E8').Value = '  +0.01%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.544'
$ws.Range('E9').Value = '  -5.96%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '35.32'
$ws.Range('E10').Value = '  -9.16%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0802'
$ws.Range('E11').Value = '  -4.11%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '7.58'
$ws.Range('E12').Value = '  -6.44%  '
$ws.Range('D13').Value = '2.935.79'
$ws.Range('E13').Value = '  -2.02%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.107'
$ws.Range('E14').Value = '  +0.04%  '
$ws.Range('D15').Value = '2.549.27'
$ws.Range('E15').Value = '  -1.95%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.858'
$ws.Range('E16').Value = '  -5.92%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '14.02'
$ws.Range('E17').Value = '  -5.34%  '
$ws.Range('D18').Value = '42.530.22'
$ws.Range('E18').Value = '  -7.80%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.64'
$ws.Range('E19').Value = '  -1.17%  '
$ws.Range('B20').Value = 'ShibaInu'
$ws.Range('C20').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D20').Value = '0.0₃0963'
$ws.Range('E20').Value = '  -5.20%  '
$ws.Range('B21').Value = 'InternetComputer(DFINITY)'
$ws.Range('C21').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.44'
$ws.Range('E21').Value = '  -2.76%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '72.29'
$ws.Range('E22').Value = '  +0.46%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '257.03'
$ws.Range('E23').Value = '  -7.46%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.87'
$ws.Range('E24').Value = '  -5.84%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '29.20'
$ws.Range('E25').Value = '  -2.06%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.10'
$ws.Range('E26').Value = '  -5.28%  '
$ws.Range('E27').Value = '  +0.19%  '
$ws.Range('E28').Value = '  -7.41%  '
$ws.Range('E29').Value = '  -4.39%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '35.84'
$ws.Range('E30').Value = '  -5.20%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '5.85'
$ws.Range('E31').Value = '  -5.80%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '150.30'
$ws.Range('E32').Value = '  -3.70%  '
$ws.Range('E33').Value = '  -2.07%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '3.35'
$ws.Range('E34').Value = '  -6.52%  '
$ws.Range('E35').Value = '  -3.10%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.0786'
$ws.Range('E36').Value = '  -6.08%  '
$ws.Range('E37').Value = '  -7.98%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '24.12'
$ws.Range('E38').Value = '  +4.86%  '
$ws.Range('E39').Value = '  -3.70%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '15.58'
$ws.Range('E40').Value = '  -0.69%  '
$ws.Range('E41').Value = '  -5.38%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0307'
$ws.Range('E42').Value = '  -6.90%  '
$ws.Range('B43').Value = 'RenderToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.78'
$ws.Range('E43').Value = '  -4.50%  '
$ws.Range('B44').Value = 'Maker'
$ws.Range('C44').Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range('D44').Value = '2.057.87'
$ws.Range('E44').Value = '  -1.83%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.999'
$ws.Range('E45').Value = '  +0.04%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '83.85'
$ws.Range('E46').Value = '  -12.07%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '1.58'
$ws.Range('E47').Value = '  +2.53%  '
$ws.Range('D48').Value = '2.790.36'
$ws.Range('E48').Value = '  -2.11%  '
$ws.Range('E49').Value = '  -9.12%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.69'
$ws.Range('E50').Value = '  -3.37%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '102.79'
$ws.Range('E51').Value = '  -5.28%  '
